# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text looks like a plain number (single decimal point) need to be
# forced to Text format first, otherwise Excel will auto-convert them to a numeric value
# and introduce floating point rounding artifacts.
$textForceCells = @("D5", "D6", "D11", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D41", "D44", "D46", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.497.39"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "2.651.80"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "596.02"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "167.39"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "2.652.80"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "5.27"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "28.07"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "3.138.05"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("D17").Value = "67.578.75"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "2.652.19"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "12.15"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +6.84%  "
$ws.Range("D21").Value = "363.25"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "4.80"
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  +9.55%  "
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "71.13"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "558.19"
$ws.Range("E31").Value = "  -5.29%  "
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").Value = "1.39"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").Value = "157.69"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "40.38"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("D48").Value = "0.595"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "154.27"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "3.87"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("E51").Value = "  -2.99%  "

# Restore the default cell style (no custom number format) now that the text values are
# safely stored as strings, so the cells keep matching the original workbook's formatting.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
